# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - timestamp string updated
#  - Suiza (row 13), Austria (row 17), Dinamarca (row 27), Rumania (row 32),
#    Libano (row 82) and Sri Lanka (row 113) get refreshed counters
#  - Senegal overtakes Montenegro in total cases, so the two swap ranking
#    positions (Senegal -> row 105 with its new data, Montenegro -> row 106
#    keeping its previous data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 12:22"

# Row 13 - Suiza
$ws.Cells.Item(13, 2).Value = 21282
$ws.Cells.Item(13, 3).Value = 182
$ws.Cells.Item(13, 5).Value = 13250
$ws.Cells.Item(13, 7).Value = 19
$ws.Cells.Item(13, 8).Value = 734

# Row 17 - Austria
$ws.Cells.Item(17, 2).Value = 12083
$ws.Cells.Item(17, 3).Value = 32
$ws.Cells.Item(17, 5).Value = 8400

# Row 27 - Dinamarca
$ws.Cells.Item(27, 6).Value = 144

# Row 32 - Rumania
$ws.Cells.Item(32, 2).Value = 4057
$ws.Cells.Item(32, 3).Value = 193
$ws.Cells.Item(32, 4).Value = 406
$ws.Cells.Item(32, 5).Value = 3494
$ws.Cells.Item(32, 6).Value = 179

# Row 82 - Libano
$ws.Cells.Item(82, 2).Value = 541
$ws.Cells.Item(82, 3).Value = 14
$ws.Cells.Item(82, 5).Value = 467
$ws.Cells.Item(82, 6).Value = 27
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = 19

# Row 105 - now Senegal (updated numbers, overtakes Montenegro)
$ws.Cells.Item(105, 1).Value = "Senegal"
$ws.Cells.Item(105, 2).Value = 226
$ws.Cells.Item(105, 3).Value = 4
$ws.Cells.Item(105, 4).Value = 92
$ws.Cells.Item(105, 5).Value = 132
$ws.Cells.Item(105, 6).Value = 1

# Row 106 - now Montenegro (keeps its previous numbers)
$ws.Cells.Item(106, 1).Value = "Montenegro"
$ws.Cells.Item(106, 2).Value = 223
$ws.Cells.Item(106, 3).Value = 9
$ws.Cells.Item(106, 4).Value = 1
$ws.Cells.Item(106, 5).Value = 220
$ws.Cells.Item(106, 6).Value = 4

# Row 113 - Sri Lanka
$ws.Cells.Item(113, 4).Value = 34
$ws.Cells.Item(113, 5).Value = 137
